# Add files via upload
# Populate the newly-added "Items Due" entries (column F) on the Burndown
# sheet for rows 11-14, which were previously blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown")

$ws.Activate()

$ws.Range("F11").Value = "17, 18, 19"
$ws.Range("F12").Value = 29
$ws.Range("F13").Value = "25, 26, 30, 31, 32, 33, 35, 34, 36, 37"
$ws.Range("F14").Value = 38

$ws.Range("D17").Select()
